$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: update rows 505-533 in place (D, K, L, M, N, O, P, S columns)
$ws.Range("D505").Value2 = 44516
$ws.Range("M505").Value2 = 160
$ws.Range("N505").Value2 = 16000
$ws.Range("O505").Value2 = 16000
$ws.Range("P505").Value2 = 16000
$ws.Range("S505").Value2 = 800
$ws.Range("D506").Value2 = 44516
$ws.Range("M506").Value2 = 280
$ws.Range("N506").Value2 = 17000
$ws.Range("O506").Value2 = 17000
$ws.Range("P506").Value2 = 17000
$ws.Range("S506").Value2 = 850
$ws.Range("D507").Value2 = 44516
$ws.Range("K507").Value2 = "Sin especificar"
$ws.Range("L507").Value2 = "Primera Pintón"
$ws.Range("M507").Value2 = 240
$ws.Range("D508").Value2 = 44270
$ws.Range("M508").Value2 = 280
$ws.Range("N508").Value2 = 12000
$ws.Range("O508").Value2 = 12500
$ws.Range("P508").Value2 = 12286
$ws.Range("S508").Value2 = 614
$ws.Range("D509").Value2 = 44270
$ws.Range("M509").Value2 = 320
$ws.Range("N509").Value2 = 13000
$ws.Range("O509").Value2 = 13000
$ws.Range("P509").Value2 = 13000
$ws.Range("S509").Value2 = 650
$ws.Range("D510").Value2 = 44295
$ws.Range("K510").Value2 = "Barraganete"
$ws.Range("L510").Value2 = "Verde"
$ws.Range("M510").Value2 = 80
$ws.Range("N510").Value2 = 18000
$ws.Range("O510").Value2 = 18000
$ws.Range("P510").Value2 = 18000
$ws.Range("S510").Value2 = 900
$ws.Range("D511").Value2 = 44295
$ws.Range("L511").Value2 = "Maduro"
$ws.Range("M511").Value2 = 160
$ws.Range("N511").Value2 = 11000
$ws.Range("O511").Value2 = 11000
$ws.Range("P511").Value2 = 11000
$ws.Range("S511").Value2 = 550
$ws.Range("D512").Value2 = 44295
$ws.Range("L512").Value2 = "Pintón"
$ws.Range("M512").Value2 = 520
$ws.Range("N512").Value2 = 12000
$ws.Range("O512").Value2 = 12500
$ws.Range("P512").Value2 = 12192
$ws.Range("S512").Value2 = 610
$ws.Range("D513").Value2 = 44217
$ws.Range("L513").Value2 = "Maduro"
$ws.Range("M513").Value2 = 200
$ws.Range("N513").Value2 = 11000
$ws.Range("O513").Value2 = 11000
$ws.Range("P513").Value2 = 11000
$ws.Range("S513").Value2 = 550
$ws.Range("D514").Value2 = 44217
$ws.Range("L514").Value2 = "Pintón"
$ws.Range("M514").Value2 = 280
$ws.Range("N514").Value2 = 12000
$ws.Range("O514").Value2 = 12000
$ws.Range("P514").Value2 = 12000
$ws.Range("S514").Value2 = 600
$ws.Range("D515").Value2 = 44509
$ws.Range("L515").Value2 = "Maduro"
$ws.Range("M515").Value2 = 160
$ws.Range("N515").Value2 = 17000
$ws.Range("O515").Value2 = 17000
$ws.Range("P515").Value2 = 17000
$ws.Range("S515").Value2 = 850
$ws.Range("D516").Value2 = 44509
$ws.Range("L516").Value2 = "Pintón"
$ws.Range("M516").Value2 = 240
$ws.Range("N516").Value2 = 18000
$ws.Range("O516").Value2 = 18000
$ws.Range("P516").Value2 = 18000
$ws.Range("S516").Value2 = 900
$ws.Range("D517").Value2 = 44509
$ws.Range("L517").Value2 = "Primera Pintón"
$ws.Range("M517").Value2 = 240
$ws.Range("N517").Value2 = 20000
$ws.Range("O517").Value2 = 20000
$ws.Range("P517").Value2 = 20000
$ws.Range("S517").Value2 = 1000
$ws.Range("D518").Value2 = 44421
$ws.Range("L518").Value2 = "Pintón"
$ws.Range("M518").Value2 = 240
$ws.Range("N518").Value2 = 13000
$ws.Range("O518").Value2 = 13000
$ws.Range("P518").Value2 = 13000
$ws.Range("S518").Value2 = 650
$ws.Range("D519").Value2 = 44421
$ws.Range("L519").Value2 = "Primera Pintón"
$ws.Range("M519").Value2 = 200
$ws.Range("N519").Value2 = 15000
$ws.Range("O519").Value2 = 15000
$ws.Range("P519").Value2 = 15000
$ws.Range("S519").Value2 = 750
$ws.Range("D520").Value2 = 44421
$ws.Range("L520").Value2 = "Verde"
$ws.Range("M520").Value2 = 120
$ws.Range("N520").Value2 = 13000
$ws.Range("O520").Value2 = 13000
$ws.Range("P520").Value2 = 13000
$ws.Range("S520").Value2 = 650
$ws.Range("D521").Value2 = 44383
$ws.Range("M521").Value2 = 160
$ws.Range("N521").Value2 = 9000
$ws.Range("O521").Value2 = 9000
$ws.Range("P521").Value2 = 9000
$ws.Range("S521").Value2 = 450
$ws.Range("D522").Value2 = 44383
$ws.Range("M522").Value2 = 240
$ws.Range("N522").Value2 = 10000
$ws.Range("O522").Value2 = 10000
$ws.Range("P522").Value2 = 10000
$ws.Range("S522").Value2 = 500
$ws.Range("D523").Value2 = 44383
$ws.Range("L523").Value2 = "Primera Pintón"
$ws.Range("M523").Value2 = 200
$ws.Range("N523").Value2 = 12000
$ws.Range("O523").Value2 = 12000
$ws.Range("P523").Value2 = 12000
$ws.Range("S523").Value2 = 600
$ws.Range("D524").Value2 = 44244
$ws.Range("L524").Value2 = "Maduro"
$ws.Range("M524").Value2 = 120
$ws.Range("N524").Value2 = 10000
$ws.Range("O524").Value2 = 10000
$ws.Range("P524").Value2 = 10000
$ws.Range("S524").Value2 = 500
$ws.Range("D525").Value2 = 44244
$ws.Range("L525").Value2 = "Pintón"
$ws.Range("M525").Value2 = 320
$ws.Range("N525").Value2 = 11000
$ws.Range("O525").Value2 = 11000
$ws.Range("P525").Value2 = 11000
$ws.Range("S525").Value2 = 550
$ws.Range("D526").Value2 = 44307
$ws.Range("M526").Value2 = 120
$ws.Range("D527").Value2 = 44307
$ws.Range("M527").Value2 = 320
$ws.Range("O527").Value2 = 12000
$ws.Range("P527").Value2 = 12000
$ws.Range("S527").Value2 = 600
$ws.Range("D528").Value2 = 44307
$ws.Range("L528").Value2 = "Primera Pintón"
$ws.Range("M528").Value2 = 160
$ws.Range("N528").Value2 = 13000
$ws.Range("O528").Value2 = 13000
$ws.Range("P528").Value2 = 13000
$ws.Range("S528").Value2 = 650
$ws.Range("D529").Value2 = 44273
$ws.Range("L529").Value2 = "Maduro"
$ws.Range("M529").Value2 = 160
$ws.Range("N529").Value2 = 11000
$ws.Range("O529").Value2 = 11000
$ws.Range("P529").Value2 = 11000
$ws.Range("S529").Value2 = 550
$ws.Range("D530").Value2 = 44273
$ws.Range("L530").Value2 = "Pintón"
$ws.Range("M530").Value2 = 480
$ws.Range("N530").Value2 = 12000
$ws.Range("O530").Value2 = 12500
$ws.Range("P530").Value2 = 12208
$ws.Range("S530").Value2 = 610
$ws.Range("D531").Value2 = 44433
$ws.Range("D532").Value2 = 44433
$ws.Range("M532").Value2 = 320
$ws.Range("N532").Value2 = 13000
$ws.Range("P532").Value2 = 13000
$ws.Range("S532").Value2 = 650
$ws.Range("D533").Value2 = 44433
$ws.Range("M533").Value2 = 200
$ws.Range("N533").Value2 = 15000
$ws.Range("O533").Value2 = 15000
$ws.Range("P533").Value2 = 15000
$ws.Range("S533").Value2 = 750

# Step 2: insert 3 rows before row 534 (old row 534 shifts to row 537)
$ws.Rows.Item(534).Insert()
$ws.Rows.Item(534).Insert()
$ws.Rows.Item(534).Insert()

# Step 3: populate new rows 534-536, and restore full row 537 (shifted content lost all but D cell)
# row 534
$ws.Range("A534").Value2 = 3
$ws.Range("B534").Value2 = "Femacal de La Calera"
$ws.Range("C534").Value2 = "Coquimbo"
$ws.Range("D534").Value2 = 44302
$ws.Range("E534").Value2 = 5
$ws.Range("F534").Value2 = "Fruta"
$ws.Range("G534").Value2 = 100108
$ws.Range("H534").Value2 = "Tropicales y subtropicales"
$ws.Range("I534").Value2 = 100108006
$ws.Range("J534").Value2 = "Plátano"
$ws.Range("K534").Value2 = "Sin especificar"
$ws.Range("L534").Value2 = "Maduro"
$ws.Range("M534").Value2 = 120
$ws.Range("N534").Value2 = 12000
$ws.Range("O534").Value2 = 12000
$ws.Range("P534").Value2 = 12000
$ws.Range("Q534").Value2 = "$/caja 20 kilos"
$ws.Range("R534").Value2 = "Ecuador"
$ws.Range("S534").Value2 = 600
$ws.Range("T534").Value2 = 20
# row 535
$ws.Range("A535").Value2 = 3
$ws.Range("B535").Value2 = "Femacal de La Calera"
$ws.Range("C535").Value2 = "Coquimbo"
$ws.Range("D535").Value2 = 44302
$ws.Range("E535").Value2 = 5
$ws.Range("F535").Value2 = "Fruta"
$ws.Range("G535").Value2 = 100108
$ws.Range("H535").Value2 = "Tropicales y subtropicales"
$ws.Range("I535").Value2 = 100108006
$ws.Range("J535").Value2 = "Plátano"
$ws.Range("K535").Value2 = "Sin especificar"
$ws.Range("L535").Value2 = "Pintón"
$ws.Range("M535").Value2 = 400
$ws.Range("N535").Value2 = 12500
$ws.Range("O535").Value2 = 13000
$ws.Range("P535").Value2 = 12800
$ws.Range("Q535").Value2 = "$/caja 20 kilos"
$ws.Range("R535").Value2 = "Ecuador"
$ws.Range("S535").Value2 = 640
$ws.Range("T535").Value2 = 20
# row 536
$ws.Range("A536").Value2 = 3
$ws.Range("B536").Value2 = "Femacal de La Calera"
$ws.Range("C536").Value2 = "Coquimbo"
$ws.Range("D536").Value2 = 44302
$ws.Range("E536").Value2 = 5
$ws.Range("F536").Value2 = "Fruta"
$ws.Range("G536").Value2 = 100108
$ws.Range("H536").Value2 = "Tropicales y subtropicales"
$ws.Range("I536").Value2 = 100108006
$ws.Range("J536").Value2 = "Plátano"
$ws.Range("K536").Value2 = "Sin especificar"
$ws.Range("L536").Value2 = "Primera Pintón"
$ws.Range("M536").Value2 = 120
$ws.Range("N536").Value2 = 14000
$ws.Range("O536").Value2 = 14000
$ws.Range("P536").Value2 = 14000
$ws.Range("Q536").Value2 = "$/caja 20 kilos"
$ws.Range("R536").Value2 = "Ecuador"
$ws.Range("S536").Value2 = 700
$ws.Range("T536").Value2 = 20
# row 537
$ws.Range("A537").Value2 = 3
$ws.Range("B537").Value2 = "Femacal de La Calera"
$ws.Range("C537").Value2 = "Coquimbo"
$ws.Range("D537").Value2 = 44179
$ws.Range("E537").Value2 = 5
$ws.Range("F537").Value2 = "Fruta"
$ws.Range("G537").Value2 = 100108
$ws.Range("H537").Value2 = "Tropicales y subtropicales"
$ws.Range("I537").Value2 = 100108006
$ws.Range("J537").Value2 = "Plátano"
$ws.Range("K537").Value2 = "Sin especificar"
$ws.Range("L537").Value2 = "Pintón"
$ws.Range("M537").Value2 = 600
$ws.Range("N537").Value2 = 13000
$ws.Range("O537").Value2 = 13500
$ws.Range("P537").Value2 = 13200
$ws.Range("Q537").Value2 = "$/caja 20 kilos"
$ws.Range("R537").Value2 = "Ecuador"
$ws.Range("S537").Value2 = 660
$ws.Range("T537").Value2 = 20
